$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- List A block (rows 2-11): compare against List B ($A$14:$A$19) ---
# Row 2 gets its own (non-shared) formula
$ws.Range("C2").Formula = '=IF(COUNTIF($A$14:$A$19,A2),A2,"")'
$ws.Range("D2").Formula = '=IF(COUNTIF($A$14:$A$19,A2),"",A2)'

# Rows 3-11 share one formula each for C and D
$ws.Range("C3:C11").Formula = '=IF(COUNTIF($A$14:$A$19,A3),A3,"")'
$ws.Range("D3:D11").Formula = '=IF(COUNTIF($A$14:$A$19,A3),"",A3)'

# --- List B block (rows 14-19): compare against List A ($A$2:$A$11) ---
# Row 14 gets its own (non-shared) formula
$ws.Range("C14").Formula = '=IF(COUNTIF($A$2:$A$11,A14),A14,"")'
$ws.Range("D14").Formula = '=IF(COUNTIF($A$2:$A$11,A14),"",A14)'

# Rows 15-19 share one formula each for C and D
$ws.Range("C15:C19").Formula = '=IF(COUNTIF($A$2:$A$11,A15),A15,"")'
$ws.Range("D15:D19").Formula = '=IF(COUNTIF($A$2:$A$11,A15),"",A15)'

# --- Selection moved to G24 ---
[void]$ws.Range("G24").Select()
